# 09-002-B Field Data.xlsx - "Revised 2015 data and results"
#
# 1. Rename the column-E header from "Common" to "SPECIES_CODE".
# 2. Sort the data rows (A2:H31) by the ORIGIN column (G), ascending.
#    At this point the substrate (rock / bare ground) rows still have a
#    blank ORIGIN, so Excel's ascending sort puts them last, after the
#    "E", "N" and "U" groups - this reproduces the row order seen in the
#    saved workbook.
# 3. Now that the substrate rows have sorted to the bottom, fill their
#    (still blank) ORIGIN cells in with "S".
# 4. Move the active selection to I31, matching the saved cursor
#    position recorded in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rename -------------------------------------------------
$ws.Range("E1").Value = "SPECIES_CODE"

# --- 2. Sort A2:H31 by column G (ORIGIN) ascending ---------------------
# Blank ORIGIN cells (the rock/bare ground substrate rows) sort to the
# bottom of the range, below "E", "N" and "U".
$sortRange = $ws.Range("A2:H31")
$sortKey = $ws.Range("G2:G31")
[void]$sortRange.Sort($sortKey)

# --- 3. Tag the now-trailing substrate rows with ORIGIN = "S" ----------
for ($r = 2; $r -le 31; $r++) {
    $originCell = $ws.Cells.Item($r, 7)
    if ([string]::IsNullOrEmpty($originCell.Value2)) {
        $originCell.Value = "S"
    }
}

# --- 4. Restore the saved selection ------------------------------------
[void]$ws.Range("I31").Select()
